# Append one new log row to each of the four lifter-log sheets, matching
# the existing table layout (columns A-I) and per-sheet formatting
# conventions already present in the workbook.

$wb = $excel.ActiveWorkbook

$dateFmt = "YYYY-MM-DD HH:MM:SS"

function Add-LogRow {
    param($ws, $row, $timeSerial, $b, $c, $d, $e, $f, $g, $h, $i)

    $ws.Cells.Item($row, 1).Value = $timeSerial
    $ws.Cells.Item($row, 1).NumberFormat = $dateFmt

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e

    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

function Add-LogRowTextId {
    param($ws, $row, $timeSerial, $b, $c, $d, $e, $f, $gText, $h, $i)

    $ws.Cells.Item($row, 1).Value = $timeSerial
    $ws.Cells.Item($row, 1).NumberFormat = $dateFmt

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e

    $ws.Cells.Item($row, 6).Value = $f

    # This particular sheet already stores the large "ID_DEC" value as a
    # literal digit string rather than a floating-point number (see the
    # existing rows above). Force text entry so the full-precision digit
    # string is preserved exactly as logged, then drop back to the
    # default style (no quote-prefix / number-format clutter) so it
    # matches how the rest of the column is stored.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = $gText
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

$idDecNumber = "5.68631262647114e+23" -as [double]
$idDecText = "568631262647113771663628"

# ROW50-FE-LIFTER
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-LogRow $ws1 61 45754.72118694444 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x5a" "0xe" 400 $idDecNumber 346 14

# ROW50-MID-LIFTER (this sheet keeps the ID_DEC column as literal digit text)
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-LogRowTextId $ws2 63 45754.68796296296 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x5e" "0x19" 400 $idDecText 350 25

# ROW11-FE-LIFTER
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-LogRow $ws3 61 45754.75482247685 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x5a" "0x14" 400 $idDecNumber 346 20

# ROW11-MID-LIFTER
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-LogRow $ws4 61 45754.88322623842 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x62" "0x19" 400 $idDecNumber 354 25
